$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 0.1920809093455981
$ws.Range("C7").Value = 0.6594294039653735
$ws.Range("D7").Value = 1.284944221537631
$ws.Range("E7").Value = 1.133553801783414
$ws.Range("F7").Value = 1.133008249231015
$ws.Range("G7").Value = 36

# Row 8
$ws.Range("B8").Value = 0.1796775225968534
$ws.Range("C8").Value = 0.6623220035822726
$ws.Range("D8").Value = 1.39608628180706
$ws.Range("E8").Value = 1.18156095137198
$ws.Range("F8").Value = 1.184868811785385
$ws.Range("G8").Value = 35

# Row 9
$ws.Range("B9").Value = -0.02192181864440764
$ws.Range("C9").Value = 0.6050125953328321
$ws.Range("D9").Value = 0.8025268610977779
$ws.Range("E9").Value = 0.8958386356357811
$ws.Range("F9").Value = 0.9188358165951085
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = 0.2162555955008903
$ws.Range("C10").Value = 0.5414520406110113
$ws.Range("D10").Value = 0.7908809848366183
$ws.Range("E10").Value = 0.8893148963312255
$ws.Range("F10").Value = 0.8978441090219604
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = 0.1228779891785326
$ws.Range("C11").Value = 0.3629836300342202
$ws.Range("D11").Value = 0.2439168504478103
$ws.Range("E11").Value = 0.493879388563453
$ws.Range("F11").Value = 0.5348105391435956
$ws.Range("G11").Value = 5
